$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'327.34"
$ws.Range("E2").Value = "'3.09%"
$ws.Range("G2").Value = "'22"
$ws.Range("D3").Value = "'40.28"
$ws.Range("E3").Value = "'5.97%"
$ws.Range("G3").Value = "'22"
$ws.Range("D4").Value = "'5.796"
$ws.Range("E4").Value = "'11.95%"
$ws.Range("G4").Value = "'22"
$ws.Range("D5").Value = "'0.08007"
$ws.Range("E5").Value = "'0.37%"
$ws.Range("G5").Value = "'22"
$ws.Range("D6").Value = "'4.595"
$ws.Range("E6").Value = "'2.75%"
$ws.Range("G6").Value = "'22"
$ws.Range("D7").Value = "'8.738"
$ws.Range("E7").Value = "'2.80%"
$ws.Range("G7").Value = "'22"
$ws.Range("D8").Value = "'1.939"
$ws.Range("E8").Value = "'0.16%"
$ws.Range("G8").Value = "'22"
$ws.Range("E9").Value = "'-1.48%"
$ws.Range("G9").Value = "'22"
$ws.Range("D10").Value = "'0.9450"
$ws.Range("E10").Value = "'0.30%"
$ws.Range("G10").Value = "'22"
$ws.Range("D11").Value = "'0.1257"
$ws.Range("E11").Value = "'0.84%"
$ws.Range("G11").Value = "'22"
$ws.Range("D12").Value = "'0.1964"
$ws.Range("E12").Value = "'1.63%"
$ws.Range("G12").Value = "'22"
$ws.Range("D13").Value = "'9.011"
$ws.Range("E13").Value = "'37.60%"
$ws.Range("G13").Value = "'22"
$ws.Range("D14").Value = "'0.09151"
$ws.Range("E14").Value = "'0.81%"
$ws.Range("G14").Value = "'22"
$ws.Range("D15").Value = "'0.03504"
$ws.Range("E15").Value = "'2.88%"
$ws.Range("G15").Value = "'22"
$ws.Range("D16").Value = "'0.09635"
$ws.Range("E16").Value = "'0.98%"
$ws.Range("G16").Value = "'22"
$ws.Range("D17").Value = "'0.001314"
$ws.Range("E17").Value = "'-3.53%"
$ws.Range("G17").Value = "'22"
$ws.Range("D18").Value = "'0.006445"
$ws.Range("E18").Value = "'8.14%"
$ws.Range("G18").Value = "'22"
$ws.Range("D19").Value = "'3.364"
$ws.Range("E19").Value = "'-1.80%"
$ws.Range("G19").Value = "'22"
$ws.Range("D20").Value = "'0.3527"
$ws.Range("E20").Value = "'0.34%"
$ws.Range("G20").Value = "'22"
$ws.Range("D21").Value = "'0.1422"
$ws.Range("E21").Value = "'8.85%"
$ws.Range("G21").Value = "'22"
$ws.Range("D22").Value = "'0.2430"
$ws.Range("E22").Value = "'5.48%"
$ws.Range("G22").Value = "'22"
$ws.Range("D23").Value = "'0.04417"
$ws.Range("E23").Value = "'1.01%"
$ws.Range("G23").Value = "'22"
$ws.Range("D24").Value = "'0.001266"
$ws.Range("E24").Value = "'3.32%"
$ws.Range("G24").Value = "'22"
$ws.Range("D25").Value = "'0.004319"
$ws.Range("E25").Value = "'-2.08%"
$ws.Range("G25").Value = "'22"
$ws.Range("D26").Value = "'0.0001145"
$ws.Range("E26").Value = "'-13.48%"
$ws.Range("G26").Value = "'22"
$ws.Range("E27").Value = "'0.36%"
$ws.Range("G27").Value = "'22"
$ws.Range("G28").Value = "'22"
$ws.Range("G29").Value = "'22"
$ws.Range("G30").Value = "'22"
$ws.Range("G31").Value = "'22"
$ws.Range("G32").Value = "'22"
$ws.Range("G33").Value = "'22"
$ws.Range("G34").Value = "'22"
$ws.Range("G35").Value = "'22"
$ws.Range("G36").Value = "'22"
$ws.Range("G37").Value = "'22"
$ws.Range("G38").Value = "'22"
$ws.Range("D39").Value = "'0.02414"
$ws.Range("E39").Value = "'-0.04%"
$ws.Range("G39").Value = "'22"
$ws.Range("D40").Value = "'0.05229"
$ws.Range("E40").Value = "'0.97%"
$ws.Range("G40").Value = "'22"
$ws.Range("D41").Value = "'0.007470"
$ws.Range("E41").Value = "'0.54%"
$ws.Range("G41").Value = "'22"
$ws.Range("D42").Value = "'0.1420"
$ws.Range("E42").Value = "'1.45%"
$ws.Range("G42").Value = "'22"
$ws.Range("D43").Value = "'0.008716"
$ws.Range("E43").Value = "'3.27%"
$ws.Range("G43").Value = "'22"
$ws.Range("D44").Value = "'0.002111"
$ws.Range("E44").Value = "'4.01%"
$ws.Range("G44").Value = "'22"
$ws.Range("D45").Value = "'0.01095"
$ws.Range("E45").Value = "'25.41%"
$ws.Range("G45").Value = "'22"
$ws.Range("D46").Value = "'0.00006952"
$ws.Range("E46").Value = "'7.17%"
$ws.Range("G46").Value = "'22"
$ws.Range("D47").Value = "'0.00000000757"
$ws.Range("E47").Value = "'1.41%"
$ws.Range("G47").Value = "'22"
$ws.Range("D48").Value = "'0.003166"
$ws.Range("E48").Value = "'10.96%"
$ws.Range("G48").Value = "'22"
$ws.Range("D49").Value = "'0.001434"
$ws.Range("E49").Value = "'-14.89%"
$ws.Range("G49").Value = "'22"
$ws.Range("D50").Value = "'0.00002120"
$ws.Range("E50").Value = "'1.41%"
$ws.Range("G50").Value = "'22"
$ws.Range("D51").Value = "'0.0002019"
$ws.Range("E51").Value = "'1.41%"
$ws.Range("G51").Value = "'22"
